# Apply weekly data-refresh edit:
#  - Row 208 gets a brand-new data point (D=44694, J=75).
#  - Rows 209..244 each "inherit" the previous row's old D/J/K/L/M/O/P
#    values (i.e. every row's data shifts down by one row).
#  - A new row 245 is appended holding the data that used to live in
#    row 244 (the last row before the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 208
$lastRow  = 244
$newRow   = 245

# Columns that participate in the shift.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# 1) Snapshot all the "before" values for the affected columns/rows so the
#    in-place writes below don't clobber data we still need to read.
#    NOTE: use .Value2 for reads - .Value returns a non-scalar wrapper in
#    this host for some property-get paths.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Row 208 gets genuinely new values for D (fecha) and J (volumen); the
#    rest of the row is untouched.
$ws.Range("D208").Value = 44694
$ws.Range("J208").Value = 75

# 3) Rows 209..244: each row takes on the previous row's OLD values.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $prev = $snapshot[$r - 1]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $prev[$c]
    }
}

# 4) Append new row 245 = a full copy of the OLD row 244 (all columns A..R).
$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($c in $allCols) {
    $ws.Range("$c$newRow").Value = $ws.Range("$c$lastRow").Value2
}
# D/J/K/L/M/O/P on the new row must hold row 244's ORIGINAL values (captured
# in the snapshot before the shift overwrote row 244).
$orig244 = $snapshot[$lastRow]
foreach ($c in $cols) {
    $ws.Range("$c$newRow").Value = $orig244[$c]
}

# Keep the date-formatted column's number format consistent with the rest
# of column D (copying the whole .Style object doesn't stick in this host,
# so copy .NumberFormat explicitly instead).
$ws.Range("D$newRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat
